$wb = $excel.ActiveWorkbook

# ALC row 70: Consecrating Congregation | Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1474
$ws.Range("I70").Value = 1407.25
$ws.Range("J70").Value = 1533.3334
$ws.Range("K70").Value = 4221.75
$ws.Range("L70").Value = 4600.0002
$ws.Range("M70").Value = -3951.75
$ws.Range("N70").Value = -5140.0002

# ALC row 73: Curbing the Contagion (L) | Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1474
$ws.Range("I73").Value = 1407.25
$ws.Range("J73").Value = 1533.3334
$ws.Range("K73").Value = 4221.75
$ws.Range("L73").Value = 4600.0002
$ws.Range("M73").Value = -3285.75
$ws.Range("N73").Value = -6472.0002

# ALC row 106: Making Your Mark | Enchanted Palladium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 11907952
$ws.Range("I106").Value = 47620936
$ws.Range("K106").Value = 47620936
$ws.Range("M106").Value = -47620305

# ALC row 129: Practical Command | Commanding Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 176409.8
$ws.Range("J129").Value = 189699.27
$ws.Range("L129").Value = 569097.8099999999
$ws.Range("N129").Value = -579097.8099999999

# ALC row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2677.6
$ws.Range("I141").Value = 2295.3635
$ws.Range("J141").Value = 3728.75
$ws.Range("K141").Value = 6886.0905
$ws.Range("L141").Value = 11186.25
$ws.Range("M141").Value = -1706.0905
$ws.Range("N141").Value = -21546.25

# ARM row 2: Ain't Got No Ingots | Bronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 673.4
$ws.Range("I2").Value = 712.6087
$ws.Range("K2").Value = 712.6087
$ws.Range("M2").Value = -599.6087

# ARM row 97: Ore for Me | High Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 37037988
$ws.Range("I97").Value = 589.7222
$ws.Range("J97").Value = 111112780
$ws.Range("K97").Value = 589.7222
$ws.Range("L97").Value = 111112780
$ws.Range("M97").Value = -93.72220000000004
$ws.Range("N97").Value = -111113772

# ARM row 116: No Scope | Titanbronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 673.4
$ws.Range("I116").Value = 712.6087
$ws.Range("K116").Value = 712.6087
$ws.Range("M116").Value = 1581.3913

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 12811.191
$ws.Range("I132").Value = 1973.3513
$ws.Range("K132").Value = 5920.0539
$ws.Range("M132").Value = -3390.0539

# BSM row 3: Hells Bells | Bronze Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 673.4
$ws.Range("I3").Value = 712.6087
$ws.Range("K3").Value = 712.6087
$ws.Range("M3").Value = -598.6087

# CRP row 62: Splinter in the Sewers | Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 52635772
$ws.Range("I62").Value = 71432190
$ws.Range("J62").Value = 5801.2
$ws.Range("K62").Value = 71432190
$ws.Range("L62").Value = 5801.2
$ws.Range("M62").Value = -71431566
$ws.Range("N62").Value = -7049.2

# CRP row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 52635772
$ws.Range("I65").Value = 71432190
$ws.Range("J65").Value = 5801.2
$ws.Range("K65").Value = 357160950
$ws.Range("L65").Value = 29006
$ws.Range("M65").Value = -357157830
$ws.Range("N65").Value = -35246

# CRP row 86: Birch, Please | Birch Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 41702184
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877

# CRP row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 41702184
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384

# CRP row 107: Built to Last | White Oak Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1323.8
$ws.Range("I107").Value = 518.53845
$ws.Range("K107").Value = 518.53845
$ws.Range("M107").Value = 1401.46155

# CRP row 120: Kindling the Flame | Lignum Vitae Ring
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 11523.714
$ws.Range("J120").Value = 12750
$ws.Range("L120").Value = 12750
$ws.Range("N120").Value = -20008

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2054.5789
$ws.Range("I132").Value = 1411.4286
$ws.Range("K132").Value = 4234.2858
$ws.Range("M132").Value = -1704.2858

# CUL row 45: Don't Turn Up Your Nose | Sauerkraut
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 630
$ws.Range("I45").Value = 265
$ws.Range("J45").Value = 995
$ws.Range("K45").Value = 795
$ws.Range("L45").Value = 2985
$ws.Range("M45").Value = -263
$ws.Range("N45").Value = -4049

# CUL row 123: Topping Up the Pot | Zurek
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 5745
$ws.Range("I123").Value = 1500
$ws.Range("J123").Value = 9990
$ws.Range("K123").Value = 4500
$ws.Range("L123").Value = 29970
$ws.Range("M123").Value = -2050
$ws.Range("N123").Value = -34870

# CUL row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 699.66
$ws.Range("J131").Value = 699.66
$ws.Range("L131").Value = 2098.98
$ws.Range("N131").Value = -12178.98

# CUL row 132: More Mezcal | Cooking Mezcal
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 786.25
$ws.Range("I132").Value = 597.5
$ws.Range("J132").Value = 975
$ws.Range("K132").Value = 5377.5
$ws.Range("L132").Value = 8775
$ws.Range("M132").Value = -2847.5
$ws.Range("N132").Value = -13835

# GSM row 46: Burning the Midnight Oil | Fire Brand
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 11061.8
$ws.Range("J46").Value = 8827.25
$ws.Range("L46").Value = 8827.25
$ws.Range("N46").Value = -9139.25

# GSM row 52: It's My Business to Know Things | Red Coral Armillae
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20005600
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 20005600
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 20005600
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -20006118

# GSM row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1799.1666
$ws.Range("I97").Value = 1639
$ws.Range("J97").Value = 2600
$ws.Range("K97").Value = 1639
$ws.Range("L97").Value = 2600
$ws.Range("M97").Value = -1143
$ws.Range("N97").Value = -3592

# GSM row 113: Copious Crystal Cannons | Manasilver Nugget
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3593.5
$ws.Range("I113").Value = 4734.3335
$ws.Range("J113").Value = 1637.7858
$ws.Range("K113").Value = 4734.3335
$ws.Range("L113").Value = 1637.7858
$ws.Range("M113").Value = -2564.3335
$ws.Range("N113").Value = -5977.7858

# GSM row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13125.426
$ws.Range("I132").Value = 2392.9534
$ws.Range("K132").Value = 7178.860199999999
$ws.Range("M132").Value = -4648.860199999999

# LTW row 22: Skin off Their Backs | Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2713.2104
$ws.Range("I22").Value = 5112.625
$ws.Range("J22").Value = 968.1818
$ws.Range("K22").Value = 5112.625
$ws.Range("L22").Value = 968.1818
$ws.Range("M22").Value = -4817.625
$ws.Range("N22").Value = -1558.1818

# LTW row 27: Fire and Hide | Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2713.2104
$ws.Range("I27").Value = 5112.625
$ws.Range("J27").Value = 968.1818
$ws.Range("K27").Value = 5112.625
$ws.Range("L27").Value = 968.1818
$ws.Range("M27").Value = -5005.625
$ws.Range("N27").Value = -1182.1818

# LTW row 46: Supply Side Logic | Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1480.0741
$ws.Range("I46").Value = 1772.6316
$ws.Range("J46").Value = 785.25
$ws.Range("K46").Value = 1772.6316
$ws.Range("L46").Value = 785.25
$ws.Range("M46").Value = -1584.6316
$ws.Range("N46").Value = -1161.25

# LTW row 61: Spelling Me Softly | Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5227
$ws.Range("I61").Value = 2252.8572
$ws.Range("K61").Value = 2252.8572
$ws.Range("M61").Value = -2050.8572

# LTW row 113: Peace in Rest | Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5227
$ws.Range("I113").Value = 2252.8572
$ws.Range("K113").Value = 2252.8572
$ws.Range("M113").Value = -82.85719999999992

# WVR row 113: A Tender Table | Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3863745.8
$ws.Range("I113").Value = 8500
$ws.Range("K113").Value = 25500
$ws.Range("M113").Value = -23330

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1555.2963
$ws.Range("I132").Value = 1116.4445
$ws.Range("J132").Value = 2433
$ws.Range("K132").Value = 3349.3335
$ws.Range("L132").Value = 7299
$ws.Range("M132").Value = -819.3335000000002
$ws.Range("N132").Value = -12359
